$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New latency data values for rows 2-13, columns A-E
$data = @(
    @(0.8786531130876747,1690542149.590441,0.7458539710474285,1690542149.621165,0.03072357177734375),
    @(0.201715374841169,1690542154.481643,0.2189853687045643,1690542154.541678,0.06003522872924805),
    @(0.5279542566709021,1690542159.485935,0.7797755126102133,1690542159.561282,0.07534694671630859),
    @(0.6956797966963151,1690542164.496362,0.6288722801249005,1690542164.551006,0.05464363098144531),
    @(0.4628335451080051,1690542169.471658,0.4350479427603577,1690542169.541372,0.06971311569213867),
    @(0.4348792884371029,1690542174.489388,0.2487097363452859,1690542174.541356,0.05196833610534668),
    @(0.6489834815756035,1690542179.48918,0.7797306323859806,1690542179.551663,0.06248378753662109),
    @(0.3538754764930114,1690542184.469036,0.3348395081812096,1690542184.537835,0.06879878044128418),
    @(0.6022871664548919,1690542189.487218,0.3516973899359572,1690542188.857749,0.6294693946838379),
    @(0.7420584498094027,1690542194.481706,0.6957297099963836,1690542194.547652,0.06594562530517578),
    @(0.5403430749682337,1690542199.485257,0.4861224937766774,1690542199.537884,0.0526275634765625),
    @(0.633418043202033,1690542204.494275,0.2552589284259358,1690542204.527619,0.03334403038024902)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
